# The deck's bespoke "Integral" theme palette is replaced by the stock
# "Office Theme" color palette (the underlying fontScheme/fmtScheme of the
# two themes are identical - only the 12 theme colors, and the cosmetic
# theme/clrScheme "name" attributes which PowerPoint does not expose for
# scripting, differ between "Integral" and "Office Theme").
#
# MsoThemeColorSchemeIndex order used by ThemeColorScheme.Item(n):
#   1 Dark1, 2 Light1, 3 Dark2, 4 Light2,
#   5 Accent1 .. 10 Accent6, 11 Hyperlink, 12 FollowedHyperlink

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

$officeThemeColors = @(
    0,          # Dark1             - 000000
    16777215,   # Light1            - FFFFFF
    6968388,    # Dark2             - 44546A
    15132391,   # Light2            - E7E6E6
    13998939,   # Accent1           - 5B9BD5
    3243501,    # Accent2           - ED7D31
    10855845,   # Accent3           - A5A5A5
    49407,      # Accent4           - FFC000
    12874308,   # Accent5           - 4472C4
    4697456,    # Accent6           - 70AD47
    12673797,   # Hyperlink         - 0563C1
    7491477     # FollowedHyperlink - 954F72
)

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $colorScheme.Item($i).RGB = $officeThemeColors[$i - 1]
}
